# Add newly-collected dictionary words, then re-sort the whole word list
# A-Z (matching the author's "append then sort" workflow).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newWords = @(
    "관절약",
    "제거술",
    "사구체신염",
    "전임의",
    "공여의사",
    "시스타틴",
    "산양유",
    "동국화인메가플러스",
    "미지참",
    "뇨단백",
    "횡문근융해증",
    "알로에",
    "갑상선기능저하증"
)

$startRow = 214
for ($i = 0; $i -lt $newWords.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newWords[$i]
}

$lastRow = $startRow + $newWords.Length - 1

$sortRange = $ws.Range("A2:C" + $lastRow)
$sortKey = $ws.Range("A2")
$sortRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1)

Write-Host "Added $($newWords.Length) words; sorted A2:C$lastRow"
